# Drop in results from RMI script
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BEPEfCT")
$ws.Range("B2").Value = 0
